# Added Code for 1) WIP-LUI_Work Complete Flow 2) VF-System Setup_Pragati's Testcases
#
# The "Routing Master" sheet holds a single sample/staging row (row 2) that the
# automation framework (Provar) stamps with the Salesforce Engineering Item
# that was most recently created/looked-up ("Pro-PEItem-<token>" in column B,
# its corresponding Salesforce Id "a345f..." in column D). Re-running the
# flow replaces that row's Item/Id pair with the newest values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routing Master")

$ws.Range("B2").Value = "Pro-PEItem-I3FR5"
$ws.Range("D2").Value = "a345f000000uRE6AAM"
